# Auto-generated edit script: apply "Add data for 2022-09-06" updates
$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("H2").Value = 76
$ws.Range("D3").Value = 97
$ws.Range("D9").Value = 305
$ws.Range("E9").Value = 308
$ws.Range("H9").Value = 316
$ws.Range("I9").Value = 376
$ws.Range("B10").Value = 898
$ws.Range("C10").Value = 1091
$ws.Range("D10").Value = 1256
$ws.Range("E10").Value = 1550
$ws.Range("F10").Value = 1587
$ws.Range("G10").Value = 762
$ws.Range("I10").Value = 619
$ws.Range("B11").Value = 1273
$ws.Range("C11").Value = 1554
$ws.Range("D11").Value = 1731
$ws.Range("E11").Value = 2017
$ws.Range("F11").Value = 2140
$ws.Range("G11").Value = 1272
$ws.Range("H11").Value = 889
$ws.Range("I11").Value = 1242

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("H5").Value = 13
$ws.Range("E7").Value = 28
$ws.Range("C8").Value = 76
$ws.Range("F8").Value = 94
$ws.Range("C19").Value = 42
$ws.Range("B28").Value = 72
$ws.Range("C28").Value = 103
$ws.Range("D28").Value = 77
$ws.Range("I28").Value = 64
$ws.Range("D32").Value = 76
$ws.Range("F32").Value = 151
$ws.Range("C36").Value = 71
$ws.Range("H36").Value = 43
$ws.Range("E50").Value = 32
$ws.Range("B53").Value = 157
$ws.Range("D53").Value = 424
$ws.Range("E53").Value = 507
$ws.Range("F53").Value = 477
$ws.Range("H53").Value = 128
$ws.Range("I53").Value = 245
$ws.Range("I61").Value = 13
$ws.Range("E63").Value = 10
$ws.Range("H65").Value = 16
$ws.Range("C70").Value = 23
$ws.Range("E70").Value = 48
$ws.Range("G70").Value = 34
$ws.Range("B74").Value = 38
$ws.Range("D74").Value = 59
$ws.Range("E74").Value = 62
$ws.Range("H77").Value = 39
$ws.Range("C78").Value = 24
$ws.Range("C80").Value = 20
$ws.Range("E87").Value = 25
$ws.Range("D91").Value = 6
$ws.Range("B95").Value = 11
$ws.Range("E95").Value = 67
$ws.Range("F95").Value = 50
$ws.Range("E97").Value = 17
$ws.Range("B99").Value = 1273
$ws.Range("C99").Value = 1554
$ws.Range("D99").Value = 1731
$ws.Range("E99").Value = 2017
$ws.Range("F99").Value = 2140
$ws.Range("G99").Value = 1272
$ws.Range("H99").Value = 889
$ws.Range("I99").Value = 1242

# Sheet 5: Roseland
$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("H7").Value = 12
$ws.Range("H9").Value = 39

# Sheet 6: Auburn Gresham
$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("E6").Value = 14
$ws.Range("E7").Value = 28

# Sheet 8: Austin
$ws = $wb.Worksheets.Item("Austin")
$ws.Range("C7").Value = 42
$ws.Range("F7").Value = 59
$ws.Range("C8").Value = 76
$ws.Range("F8").Value = 94

# Sheet 12: Garfield Park
$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("D3").Value = 5
$ws.Range("F8").Value = 99
$ws.Range("D9").Value = 76
$ws.Range("F9").Value = 151

# Sheet 13: Chatham
$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("C7").Value = 34
$ws.Range("C8").Value = 42

# Sheet 14: Grand Crossing
$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("H7").Value = 16
$ws.Range("C8").Value = 42
$ws.Range("C9").Value = 71
$ws.Range("H9").Value = 43

# Sheet 15: Loop
$ws = $wb.Worksheets.Item("Loop")
$ws.Range("E7").Value = 48
$ws.Range("H7").Value = 43
$ws.Range("I7").Value = 62
$ws.Range("B8").Value = 121
$ws.Range("D8").Value = 366
$ws.Range("F8").Value = 422
$ws.Range("B9").Value = 157
$ws.Range("D9").Value = 424
$ws.Range("E9").Value = 507
$ws.Range("F9").Value = 477
$ws.Range("H9").Value = 128
$ws.Range("I9").Value = 245

# Sheet 16: Armour Square
$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("H5").Value = 5
$ws.Range("H7").Value = 13

# Sheet 17: Old Town
$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("C6").Value = 18
$ws.Range("E6").Value = 39
$ws.Range("G6").Value = 22
$ws.Range("C7").Value = 23
$ws.Range("E7").Value = 48
$ws.Range("G7").Value = 34

# Sheet 18: Little Italy, UIC
$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("E6").Value = 19
$ws.Range("E7").Value = 32

# Sheet 19: North Lawndale
$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("H2").Value = 4
$ws.Range("H8").Value = 16

# Sheet 21: Sheffield & DePaul
$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("C6").Value = 16
$ws.Range("C7").Value = 20

# Sheet 24: Uptown
$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("E8").Value = 17
$ws.Range("E9").Value = 25

# Sheet 25: Rush & Division
$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("C5").Value = 21
$ws.Range("C6").Value = 24

# Sheet 26: Englewood
$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("B8").Value = 44
$ws.Range("C8").Value = 64
$ws.Range("D8").Value = 43
$ws.Range("I8").Value = 30
$ws.Range("B9").Value = 72
$ws.Range("C9").Value = 103
$ws.Range("D9").Value = 77
$ws.Range("I9").Value = 64

# Sheet 31: River North
$ws = $wb.Worksheets.Item("River North")
$ws.Range("B6").Value = 35
$ws.Range("D6").Value = 49
$ws.Range("E6").Value = 58
$ws.Range("B7").Value = 38
$ws.Range("D7").Value = 59
$ws.Range("E7").Value = 62

# Sheet 37: Woodlawn
$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("E5").Value = 7
$ws.Range("E7").Value = 17

# Sheet 55: West Town
$ws = $wb.Worksheets.Item("West Town")
$ws.Range("B5").Value = 8
$ws.Range("E5").Value = 62
$ws.Range("F5").Value = 44
$ws.Range("B6").Value = 11
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 50

# Sheet 58: New City
$ws = $wb.Worksheets.Item("New City")
$ws.Range("E4").Value = 2
$ws.Range("E6").Value = 10

# Sheet 60: West Lawn
$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("D4").Value = 4
$ws.Range("D6").Value = 6
